# Fixed UpdatedVinRenewal Tests to use valid VIN data and reset to original values afterwards
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the (invalid) VIN used across the sample rows with a valid VIN
$ws.Range("A2").Value = "19XFB5F5&C"
$ws.Range("A3").Value = "19XFB5F5&C"
$ws.Range("A4").Value = "19XFB5F5&C"
$ws.Range("A5").Value = "19XFB5F5&C"

# Row 2 (Honda row): mark make/model as updated, reset restraint-related codes
$ws.Range("E2").Value = "HONDA_UPDATED"
$ws.Range("F2").Value = "TEST"
$ws.Range("AC2").Value = "X"
$ws.Range("AD2").Value = "X"
$ws.Range("AE2").Value = "X"
$ws.Range("AF2").Value = "X"

# Row 4 (Toyota row): reset make/make-text back to the generic TEST placeholder
$ws.Range("D4").Value = "TEST"
$ws.Range("E4").Value = "TEST"

# Restore the active selection to F5
$null = $ws.Range("F5").Select()

# Auto-fit column B like the other feature columns
$null = $ws.Columns.Item(2).EntireColumn.AutoFit()
